$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 (pushes "Stop wage assignment" and the
# rows below it down by one, shifting their existing hyperlinks along the
# way) and populate it with the new "Stalking No Contact Order - SNCO"
# entry, keeping the list in alphabetical order.
$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = "Stalking No Contact Order - SNCO"
$ws.Range("B39").Value = "https://www.illinoislegalaid.org/legal-information/stalking-no-contact-order-request"

# Link the url cell, then re-apply the same "Hyperlink" cell style used by
# every other url cell in column B so the new row matches its neighbors.
$ws.Hyperlinks.Add($ws.Range("B39"), "https://www.illinoislegalaid.org/legal-information/stalking-no-contact-order-request")
$ws.Range("B39").Style = "Hyperlink"

$ws.Range("B39").Select()
